$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------
$xlNone  = -4142
$xlThin  = 2      # weight
$xlContinuous = 1 # line style (solid)
$xlLeft   = -4131
$xlCenter = -4108
$xlRight  = -4152
$White    = 16777215
$Black    = 0

function Clear-Borders($rng) {
    $rng.Borders.Item(7).LineStyle  = $xlNone   # left
    $rng.Borders.Item(8).LineStyle  = $xlNone   # top
    $rng.Borders.Item(9).LineStyle  = $xlNone   # bottom
    $rng.Borders.Item(10).LineStyle = $xlNone   # right
}

function Set-TopBorder($rng) {
    $rng.Borders.Item(8).LineStyle = $xlContinuous
    $rng.Borders.Item(8).Weight = $xlThin
    $rng.Borders.Item(8).ColorIndex = 1
}

function Set-BottomBorder($rng) {
    $rng.Borders.Item(9).LineStyle = $xlContinuous
    $rng.Borders.Item(9).Weight = $xlThin
    $rng.Borders.Item(9).ColorIndex = 1
}

function Reset-Alignment($rng) {
    $rng.HorizontalAlignment = 1   # xlGeneral
    $rng.VerticalAlignment = -4107 # xlBottom (Excel default)
    $rng.WrapText = $false
}

function Set-WhiteFill($rng) {
    $rng.Interior.Pattern = 1  # xlSolid
    $rng.Interior.Color = $White
}

function Clear-Fill($rng) {
    $rng.Interior.Pattern = $xlNone
}

# ---------------------------------------------------------------------------
# 1. Insert a new row before the old row 4 so the table grows from 4 data
#    rows (title / subtitle / years / single data row / source / note) to the
#    new 6-row layout (title / subtitle / years / two data rows / source).
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Insert()

# ---------------------------------------------------------------------------
# 2. Row 1: merged title text. Range grows from a single styled cell (A1) to
#    the merged A1:I1 block; plain (no fill, no border) bold 11pt Arial,
#    centred both ways, wrapped.
# ---------------------------------------------------------------------------
$r = $ws.Range("A1:I1")
$r.Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Tsalka Municipality"
$r.Merge()
Clear-Borders $r
Clear-Fill $r
$r.HorizontalAlignment = $xlCenter
$r.VerticalAlignment = $xlCenter
$r.WrapText = $true
$r.Font.Name = "Arial"
$r.Font.Size = 11
$r.Font.Bold = $true
$r.Font.Underline = $false
$r.Font.Color = $Black
$ws.Rows.Item(1).RowHeight = 51

# ---------------------------------------------------------------------------
# 3. Row 2 - subtitle "(End of year, persons)": plain 10pt Arial, solid white
#    fill, no border, default (general/bottom, no wrap) alignment, default
#    row height.
# ---------------------------------------------------------------------------
$r = $ws.Range("A2")
$r.Value = "(End of year, persons)"
Clear-Borders $r
Reset-Alignment $r
Set-WhiteFill $r
$r.Font.Name = "Arial"
$r.Font.Size = 10
$r.Font.Bold = $false
$r.Font.Underline = $false
$r.Font.Color = $Black
$ws.Rows.Item(2).RowHeight = 14.5

# ---------------------------------------------------------------------------
# 4. Row 3 - blank header cell (A3) switches font to 11pt Sylfaen, keeps its
#    top border and no fill; the year cells (B3:I3) are untouched.
# ---------------------------------------------------------------------------
$r = $ws.Range("A3")
Reset-Alignment $r
Clear-Fill $r
Clear-Borders $r
Set-TopBorder $r
$r.Font.Name = "Sylfaen"
$r.Font.Size = 11
$r.Font.Bold = $false
$r.Font.Color = $Black

# ---------------------------------------------------------------------------
# 5. Row 4 (new) - "family with disabilities Persons " with the 2017-2024
#    series. Label: left/center/wrap, top border. Numbers: plain, no
#    alignment override, no border, custom number format.
# ---------------------------------------------------------------------------
$r = $ws.Range("A4")
$r.Value = "family with disabilities Persons "
Clear-Borders $r
Set-TopBorder $r
Set-WhiteFill $r
$r.HorizontalAlignment = $xlLeft
$r.VerticalAlignment = $xlCenter
$r.WrapText = $true
$r.Font.Name = "Arial"
$r.Font.Size = 10
$r.Font.Bold = $false
$r.Font.Color = $Black
$ws.Rows.Item(4).RowHeight = 24.75

$data4 = @(197,203,219,232,243,243,245,248)
$cols = @("B","C","D","E","F","G","H","I")
for ($i = 0; $i -lt 8; $i++) {
    $rng = $ws.Range($cols[$i] + "4")
    $rng.Value = $data4[$i]
    $rng.NumberFormat = "#\ ##0"
    Clear-Borders $rng
    Reset-Alignment $rng
    Set-WhiteFill $rng
    $rng.Font.Name = "Arial"
    $rng.Font.Size = 10
    $rng.Font.Bold = $false
    $rng.Font.Color = $Black
}

# ---------------------------------------------------------------------------
# 6. Row 5 - "disabilities Persons " with the second series; bottom border
#    on the label and on I5 (closes the table visually).
# ---------------------------------------------------------------------------
$r = $ws.Range("A5")
$r.Value = "disabilities Persons "
Clear-Borders $r
Set-BottomBorder $r
Set-WhiteFill $r
$r.HorizontalAlignment = $xlLeft
$r.VerticalAlignment = $xlCenter
$r.WrapText = $true
$r.Font.Name = "Arial"
$r.Font.Size = 10
$r.Font.Bold = $false
$r.Font.Color = $Black
$ws.Rows.Item(5).RowHeight = 21

$data5 = @(232,236,252,266,280,277,282,286)
for ($i = 0; $i -lt 8; $i++) {
    $rng = $ws.Range($cols[$i] + "5")
    $rng.Value = $data5[$i]
    $rng.NumberFormat = "#\ ##0"
    Clear-Borders $rng
    Reset-Alignment $rng
    Set-WhiteFill $rng
    $rng.Font.Name = "Arial"
    $rng.Font.Size = 10
    $rng.Font.Bold = $false
    $rng.Font.Color = $Black
}
# I5 additionally carries the bottom border that closes the table.
Set-BottomBorder $ws.Range("I5")

# ---------------------------------------------------------------------------
# 7. Row 6 - Source note (former row 5 content), now merged A6:H6, and the
#    confidentiality "Note" row is removed entirely.
# ---------------------------------------------------------------------------
$r = $ws.Range("A6:H6")
$r.Value = "Source: Ministry of Internally Displaced Persons from the Occupied Territories, Labour, Health and Social Affairs of Georgia."
$r.Merge()
Clear-Borders $r
Set-WhiteFill $r
$r.HorizontalAlignment = $xlLeft
$r.VerticalAlignment = $xlCenter
$r.WrapText = $true
$r.Font.Name = "Arial"
$r.Font.Size = 9
$r.Font.Bold = $false
$r.Font.Underline = $false
$r.Font.Color = $Black
$ws.Rows.Item(6).RowHeight = 27.75

# The old confidentiality "Note" row (now shifted down to row 7 by the
# insert above) is dropped entirely - the new layout has no such row.
$ws.Rows.Item(7).Delete()

# ---------------------------------------------------------------------------
# 8. Column width - column A shrinks; B:Q no longer need an explicit custom
#    width (closest achievable value to the 20.81640625-char target given
#    this engine's column-width quantisation).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20

# ---------------------------------------------------------------------------
# 9. Selection shown when the sheet is activated now covers the merged title.
# ---------------------------------------------------------------------------
$ws.Range("A1:I1").Select()
